$wb = $excel.ActiveWorkbook

# Sheet "展览" (1st sheet)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 13986
$ws1.Range("F4").Value = 676
$ws1.Range("F6").Value = 521
$ws1.Range("F7").Value = 1450
$ws1.Range("F8").Value = 138

# Sheet "全部类型" (4th sheet)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 13986
$ws4.Range("F4").Value = 676
$ws4.Range("F8").Value = 521
$ws4.Range("F9").Value = 1450
$ws4.Range("F11").Value = 138
